$d = $word.ActiveDocument

function Replace-Text($oldText, $newText) {
    $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
}

function Insert-ItalicParagraphAfter($matchPrefix, $newText) {
    $targetIdx = -1
    $idx = 0
    foreach ($p in $d.Paragraphs) {
        $idx = $idx + 1
        $t = $p.Range.Text
        if ($t -like $matchPrefix) {
            $targetIdx = $idx
            $p.Range.InsertParagraphAfter()
            break
        }
    }
    $newP = $d.Paragraphs.Item($targetIdx + 1)
    $newP.Range.InsertAfter($newText)
    $start = $newP.Range.Start
    $end = $start + $newText.Length
    $subRange = $d.Range($start, $end)
    $subRange.Font.Italic = $true
}

# 1. Title change
Replace-Text "LOT2058 -  Engenharia Econômica" "LOT2058 -  Análise Técnico-Econômica de Bioprocessos"

# 2. Subtitle (Heading3) change
Replace-Text "Fundamentals of Economic Engineering" "Engineering Economics"

# 3. Ativação date change
Replace-Text "Ativação: 01/01/2022" "Ativação: 01/01/2025"

# 4. Insert italic English paragraph after Objetivos paragraph
Insert-ItalicParagraphAfter "Capacitar os alunos*" "Empower students to use concepts and tools of financial mathematics for economic evaluations of engineering projects."

# 5. Insert italic English paragraph after Programa resumido paragraph
Insert-ItalicParagraphAfter "Introdução à Engenharia Econômica*" "Introduction to Economic Engineering; cost estimation; interest; cash flow; depreciation; comparison of investment alternatives; break-even point; spreadsheet applications in Economic Engineering."

# 6. Insert italic English paragraph after Programa (detailed) paragraph
Insert-ItalicParagraphAfter "1.INTRODUÇÃO À ENGENHARIA ECONÔMICA*" "1.INTRODUCTION TO ECONOMIC ENGINEERING: needs for an economic project analysis (economic engineering as a decision-making analysis tool);2.COST ESTIMATION: capital investment estimation (types of capital cost estimates; most common cost indices; methods for estimating capital investment); equipment cost estimation; production cost estimation;3.INTEREST: time variable (simple interest; compound interest); effective, nominal, and equivalent rates; equivalence relationships.4.CASH FLOW: cash flow diagram; cash flow equivalence; uniform and gradient series; cash flow preparation.5.DEPRECIATION: depreciation methods;6.COMPARISON OF INVESTMENT ALTERNATIVES: profitability criteria – Equivalent Uniform Annual Value (EUAV) method; Present Value (PV) method; Internal Rate of Return (IRR) method; Modified Internal Rate of Return (MIRR) method; equipment renewal and replacement; payback period;7.Project break-even point;8.Spreadsheet applications in Economic Engineering analysis."

# 7. Add space before NF in Critério text
Replace-Text "exercícios individuais;NF = (0,8NP + 0,2NE)" "exercícios individuais; NF = (0,8NP + 0,2NE)"
